$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to compare area across three censuses (1989 / 2002 / 2014)
# in columns B:D. The export now only needs the latest (2014) figure, and
# the "(according to the population census data)" sub-title row is no
# longer used either. Drop the obsolete columns/row so the sheet matches
# the simplified single-year layout.

# Drop the 1989 and 2002 columns - only the 2014 figures (old column D) remain.
$ws.Columns("B:C").Delete()

# Drop the now-redundant "(according to the population census data)" row.
$ws.Rows(2).Delete()

# Restore the original (taller) row height used by the simplified template.
$ws.Rows("1:6").RowHeight = 20.1
